$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add new data row (row 2): SLIDER_GB slider-battery entry with nodal pricing
$ws.Range("A2").Value = "SLIDER_GB"
$ws.Range("B2").Value = "SLIDER_GB"
$ws.Range("C2").Value = "SLIDER_GB_Owner"
$ws.Range("K2").Value = $false
$ws.Range("L2").Value = 0
$ws.Range("M2").Value = 0
$ws.Range("N2").Value = 5000
$ws.Range("O2").Value = 5000
$ws.Range("P2").Value = 5000
$ws.Range("Q2").Value = "NODAL_PRICING"
$ws.Range("R2").Value = 52
$ws.Range("S2").Value = 5

# Widen column A so the new (longer) gc_id values fit
$ws.Columns.Item(1).ColumnWidth = 8.92

# Update the active selection on the sheet
$ws.Range("C13").Select() | Out-Null

Write-Host "Edit applied"
